# Rename "Voltage_Source" sheet to "Voltage Source"
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Voltage_Source")
$ws3.Name = "Voltage Source"

# Clear the old flat-table layout and rebuild as a multi-section CYME-style report
$ws3.Cells.Clear()

$ws3.Range("A1").Value = "Positive-Sequence Voltage Source"
$ws3.Range("B1").Value = "Go to Type List"
$ws3.Range("A2").Value = "ID"
$ws3.Range("B2").Value = "Bus"
$ws3.Range("C2").Value = "Voltage (pu)"
$ws3.Range("D2").Value = "Angle (deg)"
$ws3.Range("E2").Value = "Rs (pu)"
$ws3.Range("F2").Value = "Xs (pu)"
$ws3.Range("A3").Value = "End of Positive-Sequence Voltage Source"

$ws3.Range("A5").Value = "Single-Phase Voltage Source"
$ws3.Range("B5").Value = "Go to Type List"
$ws3.Range("A6").Value = "ID"
$ws3.Range("B6").Value = "Bus1"
$ws3.Range("C6").Value = "Voltage (V)"
$ws3.Range("D6").Value = "Angle (deg)"
$ws3.Range("E6").Value = "Rs (Ohm)"
$ws3.Range("F6").Value = "Xs (Ohm)"
$ws3.Range("A7").Value = "End of Single-Phase Voltage Source"

$ws3.Range("A9").Value = "Three-Phase Voltage Source with Short-Circuit Level Data"
$ws3.Range("B9").Value = "Go to Type List"
$ws3.Range("A10").Value = "ID"
$ws3.Range("B10").Value = "Bus1"
$ws3.Range("C10").Value = "Bus2"
$ws3.Range("D10").Value = "Bus3"
$ws3.Range("E10").Value = "kV (ph-ph RMS)"
$ws3.Range("F10").Value = "Angle_a (deg)"
$ws3.Range("G10").Value = "SC1ph (MVA)"
$ws3.Range("H10").Value = "SC3ph (MVA)"
$ws3.Range("A11").Value = "SUB650WYE"
$ws3.Range("B11").Value = "650_a"
$ws3.Range("C11").Value = "650_b"
$ws3.Range("D11").Value = "650_c"
$ws3.Range("E11").Value = 4.16
$ws3.Range("F11").Value = 0
$ws3.Range("G11").Value = 200000
$ws3.Range("H11").Value = 200000
$ws3.Range("A12").Value = "End of Three-Phase Voltage Source Short-Circuit Level Data"

$ws3.Range("A14").Value = "Three-Phase Voltage Source with Sequential Data"
$ws3.Range("B14").Value = "Go to Type List"
$ws3.Range("A15").Value = "ID"
$ws3.Range("B15").Value = "Bus1"
$ws3.Range("C15").Value = "Bus2"
$ws3.Range("D15").Value = "Bus3"
$ws3.Range("E15").Value = "kV (ph-ph RMS)"
$ws3.Range("F15").Value = "Angle_a (deg)"
$ws3.Range("G15").Value = "R1 (Ohm)"
$ws3.Range("H15").Value = "X1 (Ohm)"
$ws3.Range("I15").Value = "R0 (Ohm)"
$ws3.Range("J15").Value = "X0 (Ohm)"
$ws3.Range("A16").Value = "End of Three-Phase Voltage Source Sequential Data"

# Column A width (bestFit-style autosize similar to the authored workbook)
$ws3.Columns.Item(1).ColumnWidth = 53.6

# Make "Voltage Source" the active sheet with A7 selected (matches authored view state)
$ws3.Activate() | Out-Null
$ws3.Range("A7").Select() | Out-Null
